$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Update the "estado de cuenta" detail rows (16-20): the workbook now
# carries the new period (1812) and re-sorts the existing periods into
# ascending order 1812, 1901, 1902, 1903, 1904. The "Valor Mora" amounts
# travel with their period (only the 1812/1904 rows swap values; the rest
# stay the same).
$ws.Range("E16").Value = "1812"
$ws.Range("F16").Value = 31249

$ws.Range("E17").Value = "1901"
$ws.Range("F17").Value = 31249

$ws.Range("E18").Value = "1902"
$ws.Range("F18").Value = 31249

$ws.Range("E19").Value = "1903"
$ws.Range("F19").Value = 31249

$ws.Range("E20").Value = "1904"
$ws.Range("F20").Value = 22916

# --- Reposition the logo image slightly to the left (~0.19in / 241300 EMU)
$logo = $ws.Shapes.Item(1)
$logo.Left = 680600 / 12700
$logo.Width = 975600 / 12700
$logo.Height = 612000 / 12700
